$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data to row 13: date, hours, activity description.
$ws.Range("A13").Value = 44337
$ws.Range("B13").Value = 5
$ws.Range("D13").Value = "Debugging eines Konflikts zwischen jQuery und jQueryUI"

# Move selection to D14 (next empty activity cell)
$ws.Range("D14").Select()
